$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the weekly Waargenomen (G) / Verwacht (H) observations for weeks 11-37 (rows 3-29) ---
$ws.Range("H3").Value = 3253

$ws.Range("H4").Value = 3174

$ws.Range("H5").Value = 3104

$ws.Range("H6").Value = 3024

$ws.Range("G7").Value = 4979
$ws.Range("H7").Value = 2957

$ws.Range("G8").Value = 4300
$ws.Range("H8").Value = 2915

$ws.Range("G9").Value = 3907
$ws.Range("H9").Value = 2869

$ws.Range("H10").Value = 2841

$ws.Range("H11").Value = 2821

$ws.Range("H12").Value = 2794

$ws.Range("H13").Value = 2770

$ws.Range("H14").Value = 2753

$ws.Range("H15").Value = 2735

$ws.Range("H16").Value = 2737

$ws.Range("H17").Value = 2725

$ws.Range("H18").Value = 2717

$ws.Range("H19").Value = 2723

$ws.Range("G20").Value = 2616
$ws.Range("H20").Value = 2719

$ws.Range("H21").Value = 2720

$ws.Range("H22").Value = 2707

$ws.Range("G23").Value = 2658
$ws.Range("H23").Value = 2687

$ws.Range("H24").Value = 2682

$ws.Range("G25").Value = 3203
$ws.Range("H25").Value = 2669

$ws.Range("H26").Value = 2663

$ws.Range("G27").Value = 2720
$ws.Range("H27").Value = 2667

$ws.Range("G28").Value = 2668
$ws.Range("H28").Value = 2676

$ws.Range("G29").Value = 2718
$ws.Range("H29").Value = 2698

# --- Insert a new row for week 38 right after week 37 (shifts the summary row from 31 to 32) ---
$ws.Rows("30").Insert()

$ws.Range("F30").Value = 38
$ws.Range("G30").Value = 2683
$ws.Range("H30").Value = 2729

# Give the new row the same "Oversterfte" formula as the rows above it
$ws.Range("I30").Formula = "=G30-H30"

# --- Update the selection shown when the sheet is opened ---
$ws.Range("F31").Select()
